$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The lecture table used to reserve two trailing "counter-only" rows (old
# rows 15 & 16) that only held the running A-column number with no other
# data. Drop them first so the table settles back to its real size before
# we add the new header row above it.
$ws.Range("A15:A16").EntireRow.Delete()

# Add a new row 1 for the "GitHub link for all lectures" callout; this
# pushes the whole lecture table (header + all data rows) down by one.
$ws.Rows.Item(1).Insert()

# Fill in lecture 12 (Data Visualization) and lecture 13 (Multi graph
# plotting), which now occupy the two rows that used to be counter-only.
$ws.Range("B14").Value = "Data Visualization"
$ws.Range("D15").Value = "https://www.youtube.com/watch?v=Z8b90hUig_s"
$ws.Range("B15").Value = "Multi graph plotting"
$ws.Range("D14").Value = "https://www.youtube.com/watch?v=ZtxB8DS1NEE"
$ws.Range("C14").Value = "2020-MT-L12-Visualization.pdf"
$ws.Range("C15").Value = "2020-MT-L13-MultiPlots.pdf"

$ws.Range("A1").Value = "GitHub link for all lectures"
$ws.Range("C1").Value = "https://github.com/rprustagi/EL-Programming-with-Python.git"

$ws.Range("C5").Select()
